# Auto-generated edit script applying cryptos.xlsx price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.973.26"
$ws.Range("E2").Value = "  +1.01%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.138.12"
$ws.Range("E3").Value = "  +1.59%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.69"
$ws.Range("E5").Value = "  +1.46%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.08"
$ws.Range("E6").Value = "  +1.21%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.129.14"
$ws.Range("E8").Value = "  +1.53%  "
$ws.Range("E9").Value = "  +0.48%  "
$ws.Range("E10").Value = "  +2.49%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.94"
$ws.Range("E11").Value = "  +5.67%  "
$ws.Range("E12").Value = "  +0.23%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000247"
$ws.Range("E13").Value = "  +1.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.12"
$ws.Range("E14").Value = "  -1.87%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.657.56"
$ws.Range("E15").Value = "  +1.60%  "
$ws.Range("E16").Value = "  -0.19%  "
$ws.Range("E17").Value = "  +2.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.779.38"
$ws.Range("E18").Value = "  +0.91%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.135.37"
$ws.Range("E19").Value = "  +1.51%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "465.48"
$ws.Range("E20").Value = "  +1.32%  "
$ws.Range("E21").Value = "  +1.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.730"
$ws.Range("E22").Value = "  +1.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.55"
$ws.Range("E23").Value = "  +1.59%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.39"
$ws.Range("E24").Value = "  +12.63%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.15"
$ws.Range("E25").Value = "  +1.60%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "80.87"
$ws.Range("E26").Value = "  -0.18%  "
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.84"
$ws.Range("E28").Value = "  +10.33%  "
$ws.Range("E29").Value = "  +1.68%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.36"
$ws.Range("E30").Value = "  +8.23%  "
$ws.Range("B31").Value = "FirstDigitalUSD"
$ws.Range("C31").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  +0.06%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.21"
$ws.Range("E32").Value = "  +0.41%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.114"
$ws.Range("E33").Value = "  +4.88%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "27.57"
$ws.Range("E34").Value = "  +3.65%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0862"
$ws.Range("E35").Value = "  +2.08%  "
$ws.Range("E36").Value = "  +3.20%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.14"
$ws.Range("E37").Value = "  +2.77%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.27"
$ws.Range("E38").Value = "  -1.42%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.25"
$ws.Range("E39").Value = "  -2.31%  "
$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "461.10"
$ws.Range("E40").Value = "  +5.90%  "
$ws.Range("B41").Value = "Cosmos"
$ws.Range("C41").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.37"
$ws.Range("E41").Value = "  +7.79%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "51.29"
$ws.Range("E42").Value = "  +2.12%  "
$ws.Range("E43").Value = "  +9.49%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0371"
$ws.Range("E44").Value = "  +1.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.889.12"
$ws.Range("E45").Value = "  +1.10%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.12"
$ws.Range("E46").Value = "  +10.89%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.107"
$ws.Range("E47").Value = "  -0.68%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "133.09"
$ws.Range("E48").Value = "  +7.37%  "
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("E50").Value = "  +0.60%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.22"
$ws.Range("E51").Value = "  +3.94%  "
